$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing PROD claim number (text-coerced via leading apostrophe
# so it stays a text cell / keeps its existing quotePrefix style).
$ws.Range("E2").Value = "'1120170200906"

# Add the new PREPROD environment row.
$ws.Range("A3").Value = "i-preproducciongestion.segurossura.com.ar"
$ws.Range("B3").Value = "https://i-preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do"
$ws.Range("C3").Value = "ocerutti"
$ws.Range("D3").Value = "silverarrow"

# Turn the new URL cell into a real hyperlink.
$ws.Hyperlinks.Add($ws.Range("B3"), "https://i-preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do")

$ws.Range("E3").Value = "'1120170200907"

# Match the author's last selection.
$ws.Range("N8").Select()
